# Update the "Cost Data" sheet's freight-cost-per-vehicle formulas so they
# divide by 10 (annual -> per-vehicle-decade adjustment), and apply the
# existing "$#,##0.00" (green-filled) number format to the changed cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

$ws.Range("B88").Formula = "=B54/10"
$ws.Range("C88").Formula = "=B55/10"
$ws.Range("B96").Formula = "=B87/10"

$fmt = '"$"#,##0.00'
$ws.Range("B88").NumberFormat = $fmt
$ws.Range("C88").NumberFormat = $fmt
$ws.Range("B96").NumberFormat = $fmt
